# Auto-generated edit script: adds 2026-01-31 data to violent-crime-full-year workbook
# Updates column M (year 2026 running totals) across the Citywide Totals, By Neighborhood,
# and per-neighborhood sheets, plus a couple of small 2021 (column H) corrections.

$wb = $excel.ActiveWorkbook

$sheetChanges = @(
    @{ Sheet = 'Citywide Totals'; Cells = @{ "M2" = 428; "M3" = 471; "H4" = 1771; "M4" = 134; "M5" = 29; "M6" = 363; "H7" = 26086; "M7" = 1425 } },
    @{ Sheet = 'By Neighborhood'; Cells = @{ "M2" = 12; "M7" = 42; "M8" = 95; "M9" = 12; "M19" = 50; "M20" = 48; "M27" = 20; "M29" = 71; "M33" = 52; "M34" = 9; "M37" = 66; "M42" = 48; "M48" = 17; "M60" = 12; "H63" = 322; "M63" = 5; "M65" = 30; "M68" = 4; "M71" = 9; "M72" = 10; "M76" = 17; "M77" = 15; "M78" = 22; "M79" = 34; "M83" = 31; "M85" = 69; "M87" = 4; "M88" = 17; "M89" = 21; "M90" = 14; "M91" = 20; "M94" = 19; "M95" = 19; "M96" = 13; "H101" = 26086; "M101" = 1425 } },
    @{ Sheet = 'West Ridge'; Cells = @{ "M6" = 5; "M7" = 13 } },
    @{ Sheet = 'Auburn Gresham'; Cells = @{ "M3" = 15; "M7" = 42 } },
    @{ Sheet = 'Uptown'; Cells = @{ "M4" = 5; "M7" = 21 } },
    @{ Sheet = 'South Shore'; Cells = @{ "M4" = 3; "M7" = 69 } },
    @{ Sheet = 'Austin'; Cells = @{ "M2" = 27; "M3" = 33; "M5" = 3; "M7" = 95 } },
    @{ Sheet = 'South Chicago'; Cells = @{ "M3" = 14; "M7" = 31 } },
    @{ Sheet = 'Garfield Park'; Cells = @{ "M2" = 12; "M5" = 3; "M7" = 52 } },
    @{ Sheet = 'West Pullman'; Cells = @{ "M3" = 7; "M6" = 6; "M7" = 19 } },
    @{ Sheet = 'Grand Crossing'; Cells = @{ "M6" = 17; "M7" = 66 } },
    @{ Sheet = 'New City'; Cells = @{ "M3" = 12; "M6" = 7; "M7" = 30 } },
    @{ Sheet = 'Englewood'; Cells = @{ "M2" = 23; "M3" = 23; "M4" = 6; "M6" = 18; "M7" = 71 } },
    @{ Sheet = 'Lake View'; Cells = @{ "M3" = 2; "M7" = 17 } },
    @{ Sheet = 'Chatham'; Cells = @{ "M2" = 15; "M7" = 50 } },
    @{ Sheet = 'River North'; Cells = @{ "M3" = 4; "M7" = 17 } },
    @{ Sheet = 'Humboldt Park'; Cells = @{ "M6" = 15; "M7" = 48 } },
    @{ Sheet = 'Rogers Park'; Cells = @{ "M3" = 11; "M7" = 22 } },
    @{ Sheet = 'Washington Park'; Cells = @{ "M5" = 1; "M7" = 20 } },
    @{ Sheet = 'Roseland'; Cells = @{ "M2" = 12; "M6" = 8; "M7" = 34 } },
    @{ Sheet = 'Chicago Lawn'; Cells = @{ "M2" = 20; "M3" = 11; "M4" = 6; "M7" = 48 } },
    @{ Sheet = 'Garfield Ridge'; Cells = @{ "M2" = 1; "M7" = 9 } },
    @{ Sheet = 'West Loop'; Cells = @{ "M3" = 5; "M7" = 19 } },
    @{ Sheet = 'Avalon Park'; Cells = @{ "M6" = 6; "M7" = 12 } },
    @{ Sheet = 'Albany Park'; Cells = @{ "M4" = 5; "M7" = 12 } },
    @{ Sheet = 'United Center'; Cells = @{ "M2" = 4; "M6" = 9; "M7" = 17 } },
    @{ Sheet = 'Edgewater'; Cells = @{ "M4" = 6; "M6" = 3; "M7" = 20 } },
    @{ Sheet = 'Washington Heights'; Cells = @{ "M6" = 3; "M7" = 14 } },
    @{ Sheet = 'North Park'; Cells = @{ "M6" = 1; "M7" = 4 } },
    @{ Sheet = 'Morgan Park'; Cells = @{ "M2" = 5; "M7" = 12 } },
    @{ Sheet = 'Oakland'; Cells = @{ "M4" = 1; "M7" = 9 } },
    @{ Sheet = 'Old Town'; Cells = @{ "M2" = 3; "M4" = 4; "M7" = 10 } },
    @{ Sheet = 'Riverdale'; Cells = @{ "M3" = 5; "M7" = 15 } },
    @{ Sheet = 'Ukrainian Village'; Cells = @{ "M6" = 1; "M7" = 4 } }
)

foreach ($sc in $sheetChanges) {
    $ws = $wb.Worksheets.Item($sc.Sheet)
    foreach ($cellRef in $sc.Cells.Keys) {
        $ws.Range($cellRef).Value = $sc.Cells[$cellRef]
    }
}
